$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21,2).Text = "13"

$c9 = $ws.Range("C9")
$c9.Characters(27,9).Text = "3/24/2025"
$c9.Characters(47,9).Text = "3/30/2025"

# --- Data table updates (rows 16-31) ---

# Row 16
$ws.Range("C16").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 29
$ws.Range("J16").Value = 22
$ws.Range("K16").Value = 31.818181818181
$ws.Range("L16").Value = -17.142857142857
$ws.Range("M16").Value = 31.818181818181
$ws.Range("N16").Value = -81.761006289308

# Row 17
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 5.263157894736
$ws.Range("I17").Value = 39
$ws.Range("J17").Value = 56
$ws.Range("K17").Value = -30.357142857142
$ws.Range("L17").Value = -18.75
$ws.Range("M17").Value = 8.333333333333
$ws.Range("N17").Value = -31.578947368421

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 23
$ws.Range("K18").Value = -47.727272727272
$ws.Range("L18").Value = -41.025641025641
$ws.Range("M18").Value = -43.902439024390
$ws.Range("N18").Value = -88.082901554404

# Row 19
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -21.153846153846
$ws.Range("I19").Value = 159
$ws.Range("J19").Value = 182
$ws.Range("K19").Value = -12.637362637362
$ws.Range("L19").Value = -0.625
$ws.Range("M19").Value = 25.196850393700
$ws.Range("N19").Value = -52.395209580838

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("J20").Value = 9
$ws.Range("K20").Value = -44.444444444444
$ws.Range("N20").Value = -94.252873563218

# Row 21
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -6.25
$ws.Range("F21").Value = 78
$ws.Range("G21").Value = 95
$ws.Range("H21").Value = -17.894736842105
$ws.Range("I21").Value = 259
$ws.Range("J21").Value = 314
$ws.Range("K21").Value = -17.515923566879
$ws.Range("L21").Value = -12.5
$ws.Range("M21").Value = 12.608695652173
$ws.Range("N21").Value = -68.907563025210

# Row 22
$ws.Range("D22").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 13
$ws.Range("K22").Value = 18.181818181818
$ws.Range("L22").Value = 225
$ws.Range("M22").Value = 30

# Row 23
$ws.Range("C23").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 150

# Row 24
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 27.777777777777
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = -20.689655172413
$ws.Range("I24").Value = 286
$ws.Range("J24").Value = 316
$ws.Range("K24").Value = -9.493670886075
$ws.Range("L24").Value = 16.260162601626
$ws.Range("M24").Value = 23.275862068965

# Row 25
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -16.666666666666
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 70
$ws.Range("H25").Value = -38.571428571428
$ws.Range("I25").Value = 211
$ws.Range("J25").Value = 266
$ws.Range("K25").Value = -20.676691729323
$ws.Range("L25").Value = 34.394904458598

# Row 26
$ws.Range("C26").Value = 20
$ws.Range("D26").Value = 10
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = 29.729729729729
$ws.Range("I26").Value = 109
$ws.Range("J26").Value = 95
$ws.Range("K26").Value = 14.736842105263
$ws.Range("L26").Value = 47.297297297297
$ws.Range("M26").Value = 73.015873015873

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -20

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 12
$ws.Range("K28").Value = -40
$ws.Range("L28").Value = 9.090909090909

# Row 31
$ws.Range("D31").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$ws.Range("F31").Value = 5
$ws.Range("H31").Value = 150
$ws.Range("I31").Value = 4
$ws.Range("J31").Value = 5
$ws.Range("K31").Value = 20
$ws.Range("L31").Value = 100
